$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert new front-matter sheet before the existing Sheet1, matching the
# new "Sheet2" tab that becomes the active sheet.
$newSheet = $wb.Worksheets.Add($ws1)
$newSheet.Name = "Sheet2"

$newSheet.Range("A1").Value = "New Phytologist Supporting Information"
$newSheet.Range("A1").Font.Bold = $true

$newSheet.Range("A2").Value = "Photographs as an essential biodiversity resource: drivers of gaps in the vascular plant photographic record"
$newSheet.Range("A3").Value = "Thomas Mesaglio, Hervé Sauquet, David Coleman, Elizabeth Wenk, William K Cornwell"
$newSheet.Range("A4").Value = "Accepted 8 February 2023"

$newSheet.Range("A6").Value = "Caption"
$newSheet.Range("A6").Font.Bold = $true

$newSheet.Range("A7").Value = "Definitions of the geographic abbreviations used in column headings throughout R scripts and the other supplementary files. "

$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

[void]$newSheet.Range("H22").Select()
